# "blog update and about remove card"
#
# - 2023 sheet: new blog entry image "2023_waciis.jpg" attached to the
#   existing WACIIS 2023 row.
# - 2024 sheet: new blog entry image "2024_ihci.jpg" attached to the
#   existing IHCI 2023 row, plus a custom width for column A.
# - 2022 sheet: new row for the LINCS master's-thesis blurb (June).
# - Selections / active-tab bookkeeping follow from the above edits.

$wb = $excel.ActiveWorkbook

# ---- 2024 sheet (new shared string must land first, index 54) -----------
$ws2024 = $wb.Worksheets.Item("2024")
$ws2024.Activate()
$ws2024.Range("C2").Value = "2024_ihci.jpg"
$ws2024.Columns.Item(1).ColumnWidth = 19.5
$ws2024.Range("J21").Select()

# ---- 2023 sheet (new shared string index 55) -----------------------------
$ws2023 = $wb.Worksheets.Item("2023")
$ws2023.Activate()
$ws2023.Range("C3").Value = "2023_waciis.jpg"
$ws2023.Range("I11").Select()

# ---- 2022 sheet (new shared strings index 56 + reuse "June") ------------
$ws2022 = $wb.Worksheets.Item("2022")
$ws2022.Activate()
$ws2022.Range("A4").Value = "Completed my Master's thesis at the Learning and Information in Networked Complex Systems (LINCS) Group, IIT Delhi"
$ws2022.Range("B4").Value = "June"
$ws2022.Range("C9").Select()

# ---- 2024 sheet ends active / tabSelected --------------------------------
$ws2024.Activate()
